$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.424.52'
$ws.Range('D3').Value = '3.535.79'
$ws.Range('E3').Value = '  -2.16%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '196.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '582.58'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.48%  '
$ws.Range('E7').Value = '  -2.41%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.630'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '51.79'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000287'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.26'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.04%  '
$ws.Range('D14').Value = '4.094.72'
$ws.Range('E14').Value = '  -2.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '664.46'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +12.16%  '
$ws.Range('D16').Value = '69.501.60'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '3.533.48'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.41'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.88%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.65%  '
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.968'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.80%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.27'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.30'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '105.83'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.38'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.96%  '
$ws.Range('E26').Value = '  -3.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.39'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.84'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.60%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.90'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.49%  '
$ws.Range('E33').Value = '  -5.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '61.90'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.21%  '
$ws.Range('D35').Value = '3.798.99'
$ws.Range('E35').Value = '  -2.91%  '
$ws.Range('D36').Value = '0.0₃0812'
$ws.Range('E36').Value = '  -9.48%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.68'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '504.33'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.54%  '
$ws.Range('E40').Value = '  -7.25%  '
$ws.Range('E41').Value = '  -5.04%  '
$ws.Range('E42').Value = '  +0.33%  '
$ws.Range('E43').Value = '  -7.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0454'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.02%  '
$ws.Range('E46').Value = '  +0.57%  '
$ws.Range('E47').Value = '  -3.20%  '
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('E49').Value = '  -3.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.77'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +19.74%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.72'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +62.65%  '
